$d = $word.ActiveDocument

# Locate the "footer" list-item paragraph and position the insertion point
# right after its text (before the trailing _GoBack bookmark).
$d.Content.Find.Execute("footer", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Footer", 2)

# Find the paragraph that now reads "Footer" and insert a new list paragraph
# right after it containing "Putting it all together ".
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Footer") {
        $target = $p
    }
}

$newPara = $target.Range.InsertParagraphAfter()
$afterRange = $target.Next().Range
$afterRange.Text = "Putting it all together "
